$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the two existing callback option rows:
#   enterCB -> onEnter
#   leaveCB -> onLeave
$ws.Range("A16").Value = "onEnter"
$ws.Range("A17").Value = "onLeave"

# Add a brand new row for the onChange callback option.
$ws.Range("A18").Value = "onChange"
$ws.Range("B18").Value = "Function"
$ws.Range("C18").Value = "[CodePen]()"
$ws.Range("D18").Value = "Callback function that fires when mouse moves inside of hoverTarget and a new transition is fired."

# Update the saved selection/active cell to match the authored state.
$ws.Range("C23").Select()
